$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Marks")

# Delete rows 13 through 21 (these records were removed entirely)
$ws.Range("A13:D21").EntireRow.Delete() | Out-Null

# Update the remaining data cells that changed value
$ws.Range("C6").Value = "Sub-25"
$ws.Range("D6").Value = 82

$ws.Range("D9").Value = 89

$ws.Range("C10").Value = "Sub-28"
$ws.Range("D10").Value = 81

$ws.Range("B11").Value = 5
$ws.Range("C11").Value = "Sub-22"
$ws.Range("D11").Value = 92

$ws.Range("B12").Value = 9
$ws.Range("C12").Value = "Sub-19"
$ws.Range("D12").Value = 98

# Move the active selection to G10, matching the saved view state
$ws.Range("G10").Select() | Out-Null
